$wb = $excel.ActiveWorkbook

$ws_exhibit = $wb.Worksheets.Item("展览")
$ws_exhibit.Range("F3").Value = 882
$ws_exhibit.Range("F4").Value = 882
$ws_exhibit.Range("F5").Value = 154
$ws_exhibit.Range("F6").Value = 21
$ws_exhibit.Range("F7").Value = 42
$ws_exhibit.Range("F8").Value = 2706
$ws_exhibit.Range("F11").Value = 7318
$ws_exhibit.Range("F13").Value = 7455
$ws_exhibit.Range("F16").Value = 5886
$ws_exhibit.Range("F17").Value = 3192
$ws_exhibit.Range("F18").Value = 3564
$ws_exhibit.Range("F21").Value = 263
$ws_exhibit.Range("F22").Value = 229
$ws_exhibit.Range("F23").Value = 2023
$ws_exhibit.Range("F24").Value = 98
$ws_exhibit.Range("F25").Value = 331
$ws_exhibit.Range("F28").Value = 927
$ws_exhibit.Range("F30").Value = 2532
$ws_exhibit.Range("F31").Value = 1359
$ws_exhibit.Range("F32").Value = 3064
$ws_exhibit.Range("F33").Value = 116
$ws_exhibit.Range("F35").Value = 193
$ws_exhibit.Range("F36").Value = 450
$ws_exhibit.Range("F37").Value = 1184
$ws_exhibit.Range("F38").Value = 218

$ws_show = $wb.Worksheets.Item("演出")
$ws_show.Range("F6").Value = 229

$ws_local = $wb.Worksheets.Item("本地生活")
$ws_local.Range("F2").Value = 106

$ws_all = $wb.Worksheets.Item("全部类型")
$ws_all.Range("F6").Value = 882
$ws_all.Range("F7").Value = 882
$ws_all.Range("F8").Value = 154
$ws_all.Range("F9").Value = 21
$ws_all.Range("F11").Value = 106
$ws_all.Range("F12").Value = 2706
$ws_all.Range("F14").Value = 229
$ws_all.Range("F18").Value = 7318
$ws_all.Range("F20").Value = 7455
$ws_all.Range("F23").Value = 5886
$ws_all.Range("F24").Value = 3192
$ws_all.Range("F25").Value = 3564
$ws_all.Range("F29").Value = 263
$ws_all.Range("F31").Value = 2023
$ws_all.Range("F35").Value = 331
$ws_all.Range("F37").Value = 927
$ws_all.Range("F39").Value = 2532
$ws_all.Range("F40").Value = 1359
$ws_all.Range("F42").Value = 3064
$ws_all.Range("F43").Value = 116
$ws_all.Range("F44").Value = 193
$ws_all.Range("F46").Value = 450
$ws_all.Range("F47").Value = 1184
$ws_all.Range("F48").Value = 218
